$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 659.5294
$ws.Cells.Item(33, 9).Value = 711
$ws.Cells.Item(33, 10).Value = 536
$ws.Cells.Item(33, 11).Value = 711
$ws.Cells.Item(33, 12).Value = 536
$ws.Cells.Item(33, 13).Value = -482
$ws.Cells.Item(33, 14).Value = -994
$ws.Cells.Item(70, 8).Value = 3486.3
$ws.Cells.Item(70, 9).Value = 4136.1055
$ws.Cells.Item(70, 10).Value = 2363.9092
$ws.Cells.Item(70, 11).Value = 12408.3165
$ws.Cells.Item(70, 12).Value = 7091.7276
$ws.Cells.Item(70, 13).Value = -12138.3165
$ws.Cells.Item(70, 14).Value = -7631.7276
$ws.Cells.Item(73, 8).Value = 3486.3
$ws.Cells.Item(73, 9).Value = 4136.1055
$ws.Cells.Item(73, 10).Value = 2363.9092
$ws.Cells.Item(73, 11).Value = 12408.3165
$ws.Cells.Item(73, 12).Value = 7091.7276
$ws.Cells.Item(73, 13).Value = -11472.3165
$ws.Cells.Item(73, 14).Value = -8963.7276
$ws.Cells.Item(80, 8).Value = 5035.0713
$ws.Cells.Item(80, 9).Value = 2252.08
$ws.Cells.Item(80, 10).Value = 9127.706
$ws.Cells.Item(80, 11).Value = 6756.24
$ws.Cells.Item(80, 12).Value = 27383.118
$ws.Cells.Item(80, 13).Value = -5758.24
$ws.Cells.Item(80, 14).Value = -29379.118
$ws.Cells.Item(82, 8).Value = 884.4286
$ws.Cells.Item(82, 9).Value = 340.2
$ws.Cells.Item(82, 10).Value = 2245
$ws.Cells.Item(82, 11).Value = 1020.6
$ws.Cells.Item(82, 12).Value = 6735
$ws.Cells.Item(82, 13).Value = -614.5999999999999
$ws.Cells.Item(82, 14).Value = -7547
$ws.Cells.Item(83, 8).Value = 5035.0713
$ws.Cells.Item(83, 9).Value = 2252.08
$ws.Cells.Item(83, 10).Value = 9127.706
$ws.Cells.Item(83, 11).Value = 20268.72
$ws.Cells.Item(83, 12).Value = 82149.35400000001
$ws.Cells.Item(83, 13).Value = -15276.72
$ws.Cells.Item(83, 14).Value = -92133.35400000001
$ws.Cells.Item(85, 8).Value = 884.4286
$ws.Cells.Item(85, 9).Value = 340.2
$ws.Cells.Item(85, 10).Value = 2245
$ws.Cells.Item(85, 11).Value = 1020.6
$ws.Cells.Item(85, 12).Value = 6735
$ws.Cells.Item(85, 13).Value = 383.4000000000001
$ws.Cells.Item(85, 14).Value = -9543
$ws.Cells.Item(88, 8).Value = 6420.4
$ws.Cells.Item(88, 9).Value = 1000
$ws.Cells.Item(88, 10).Value = 11840.8
$ws.Cells.Item(88, 11).Value = 1000
$ws.Cells.Item(88, 12).Value = 11840.8
$ws.Cells.Item(88, 13).Value = -594
$ws.Cells.Item(88, 14).Value = -12652.8
$ws.Cells.Item(91, 8).Value = 6420.4
$ws.Cells.Item(91, 9).Value = 1000
$ws.Cells.Item(91, 10).Value = 11840.8
$ws.Cells.Item(91, 11).Value = 1000
$ws.Cells.Item(91, 12).Value = 11840.8
$ws.Cells.Item(91, 13).Value = 404
$ws.Cells.Item(91, 14).Value = -14648.8
$ws.Cells.Item(101, 8).Value = 716.875
$ws.Cells.Item(101, 9).Value = 710
$ws.Cells.Item(101, 10).Value = 723.75
$ws.Cells.Item(101, 11).Value = 2130
$ws.Cells.Item(101, 12).Value = 2171.25
$ws.Cells.Item(101, 13).Value = -508
$ws.Cells.Item(101, 14).Value = -5415.25
$ws.Cells.Item(116, 8).Value = 4627.5293
$ws.Cells.Item(116, 9).Value = 4860.909
$ws.Cells.Item(116, 10).Value = 4199.6665
$ws.Cells.Item(116, 11).Value = 4860.909
$ws.Cells.Item(116, 12).Value = 4199.6665
$ws.Cells.Item(116, 13).Value = -1418.909
$ws.Cells.Item(116, 14).Value = -11083.6665

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 439914.03
$ws.Cells.Item(61, 9).Value = 356513.9
$ws.Cells.Item(61, 10).Value = 591996.5600000001
$ws.Cells.Item(61, 11).Value = 356513.9
$ws.Cells.Item(61, 12).Value = 591996.5600000001
$ws.Cells.Item(61, 13).Value = -356301.9
$ws.Cells.Item(61, 14).Value = -592420.5600000001
$ws.Cells.Item(74, 8).Value = 165474.16
$ws.Cells.Item(74, 9).Value = 185962.11
$ws.Cells.Item(74, 10).Value = 80370.30499999999
$ws.Cells.Item(74, 11).Value = 185962.11
$ws.Cells.Item(74, 12).Value = 80370.30499999999
$ws.Cells.Item(74, 13).Value = -185088.11
$ws.Cells.Item(74, 14).Value = -82118.30499999999
$ws.Cells.Item(77, 8).Value = 165474.16
$ws.Cells.Item(77, 9).Value = 185962.11
$ws.Cells.Item(77, 10).Value = 80370.30499999999
$ws.Cells.Item(77, 11).Value = 929810.5499999999
$ws.Cells.Item(77, 12).Value = 401851.525
$ws.Cells.Item(77, 13).Value = -925442.5499999999
$ws.Cells.Item(77, 14).Value = -410587.525
$ws.Cells.Item(88, 8).Value = 3466.6667
$ws.Cells.Item(88, 9).Value = 4916.6665
$ws.Cells.Item(88, 10).Value = 2500
$ws.Cells.Item(88, 11).Value = 4916.6665
$ws.Cells.Item(88, 12).Value = 2500
$ws.Cells.Item(88, 13).Value = -4510.6665
$ws.Cells.Item(88, 14).Value = -3312
$ws.Cells.Item(91, 8).Value = 3466.6667
$ws.Cells.Item(91, 9).Value = 4916.6665
$ws.Cells.Item(91, 10).Value = 2500
$ws.Cells.Item(91, 11).Value = 4916.6665
$ws.Cells.Item(91, 12).Value = 2500
$ws.Cells.Item(91, 13).Value = -3512.6665
$ws.Cells.Item(91, 14).Value = -5308
$ws.Cells.Item(102, 8).Value = 9678.571
$ws.Cells.Item(102, 9).Value = 1437.5
$ws.Cells.Item(102, 10).Value = 20666.666
$ws.Cells.Item(102, 11).Value = 1437.5
$ws.Cells.Item(102, 12).Value = 20666.666
$ws.Cells.Item(102, 13).Value = 184.5
$ws.Cells.Item(102, 14).Value = -23910.666
$ws.Cells.Item(136, 8).Value = 439914.03
$ws.Cells.Item(136, 9).Value = 356513.9
$ws.Cells.Item(136, 10).Value = 591996.5600000001
$ws.Cells.Item(136, 11).Value = 1069541.7
$ws.Cells.Item(136, 12).Value = 1775989.68
$ws.Cells.Item(136, 13).Value = -1066991.7
$ws.Cells.Item(136, 14).Value = -1781089.68

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 6340
$ws.Cells.Item(86, 9).Value = 13475
$ws.Cells.Item(86, 10).Value = 1583.3334
$ws.Cells.Item(86, 11).Value = 13475
$ws.Cells.Item(86, 12).Value = 1583.3334
$ws.Cells.Item(86, 13).Value = -12352
$ws.Cells.Item(86, 14).Value = -3829.3334
$ws.Cells.Item(89, 8).Value = 6340
$ws.Cells.Item(89, 9).Value = 13475
$ws.Cells.Item(89, 10).Value = 1583.3334
$ws.Cells.Item(89, 11).Value = 67375
$ws.Cells.Item(89, 12).Value = 7916.666999999999
$ws.Cells.Item(89, 13).Value = -61759
$ws.Cells.Item(89, 14).Value = -19148.667
$ws.Cells.Item(94, 8).Value = 1311.9584
$ws.Cells.Item(94, 9).Value = 1003.7857
$ws.Cells.Item(94, 10).Value = 1743.4
$ws.Cells.Item(94, 11).Value = 1003.7857
$ws.Cells.Item(94, 12).Value = 1743.4
$ws.Cells.Item(94, 13).Value = -552.7857
$ws.Cells.Item(94, 14).Value = -2645.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2756.4727
$ws.Cells.Item(31, 9).Value = 2069.889
$ws.Cells.Item(31, 11).Value = 2069.889
$ws.Cells.Item(31, 13).Value = -1774.889
$ws.Cells.Item(33, 8).Value = 2500
$ws.Cells.Item(33, 9).Value = 2500
$ws.Cells.Item(33, 11).Value = 2500
$ws.Cells.Item(33, 13).Value = -2121
$ws.Cells.Item(34, 8).Value = 2756.4727
$ws.Cells.Item(34, 9).Value = 2069.889
$ws.Cells.Item(34, 11).Value = 2069.889
$ws.Cells.Item(34, 13).Value = -1867.889
$ws.Cells.Item(62, 8).Value = 5055075.5
$ws.Cells.Item(62, 9).Value = 27779278
$ws.Cells.Item(62, 10).Value = 5253
$ws.Cells.Item(62, 11).Value = 27779278
$ws.Cells.Item(62, 12).Value = 5253
$ws.Cells.Item(62, 13).Value = -27778654
$ws.Cells.Item(62, 14).Value = -6501
$ws.Cells.Item(65, 8).Value = 5055075.5
$ws.Cells.Item(65, 9).Value = 27779278
$ws.Cells.Item(65, 10).Value = 5253
$ws.Cells.Item(65, 11).Value = 138896390
$ws.Cells.Item(65, 12).Value = 26265
$ws.Cells.Item(65, 13).Value = -138893270
$ws.Cells.Item(65, 14).Value = -32505
$ws.Cells.Item(122, 8).Value = 2081.7273
$ws.Cells.Item(122, 9).Value = 1133.3334
$ws.Cells.Item(122, 10).Value = 2437.375
$ws.Cells.Item(122, 11).Value = 3400.0002
$ws.Cells.Item(122, 12).Value = 7312.125
$ws.Cells.Item(122, 13).Value = -950.0001999999999
$ws.Cells.Item(122, 14).Value = -12212.125
$ws.Cells.Item(134, 8).Value = 1376.5518
$ws.Cells.Item(134, 9).Value = 955.04346
$ws.Cells.Item(134, 11).Value = 2865.13038
$ws.Cells.Item(134, 13).Value = -330.1303800000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 619.5357
$ws.Cells.Item(113, 10).Value = 555.3333
$ws.Cells.Item(113, 12).Value = 1665.9999
$ws.Cells.Item(113, 14).Value = -6005.9999
$ws.Cells.Item(131, 8).Value = 1407.8689
$ws.Cells.Item(131, 9).Value = 1558.8889
$ws.Cells.Item(131, 10).Value = 1381.7307
$ws.Cells.Item(131, 11).Value = 4676.6667
$ws.Cells.Item(131, 12).Value = 4145.1921
$ws.Cells.Item(131, 13).Value = 363.3333000000002
$ws.Cells.Item(131, 14).Value = -14225.1921
$ws.Cells.Item(134, 8).Value = 4231.5
$ws.Cells.Item(134, 9).Value = 4629.6665
$ws.Cells.Item(134, 10).Value = 3833.3333
$ws.Cells.Item(134, 11).Value = 13888.9995
$ws.Cells.Item(134, 12).Value = 11499.9999
$ws.Cells.Item(134, 13).Value = -8818.999500000002
$ws.Cells.Item(134, 14).Value = -21639.9999
$ws.Cells.Item(139, 8).Value = 2517.2
$ws.Cells.Item(139, 9).Value = 685.44446
$ws.Cells.Item(139, 10).Value = 4015.9092
$ws.Cells.Item(139, 11).Value = 2056.33338
$ws.Cells.Item(139, 12).Value = 12047.7276
$ws.Cells.Item(139, 13).Value = 3083.66662
$ws.Cells.Item(139, 14).Value = -22327.7276

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 4755.5117
$ws.Cells.Item(70, 9).Value = 3767.96
$ws.Cells.Item(70, 10).Value = 6127.1113
$ws.Cells.Item(70, 11).Value = 3767.96
$ws.Cells.Item(70, 12).Value = 6127.1113
$ws.Cells.Item(70, 13).Value = -3497.96
$ws.Cells.Item(70, 14).Value = -6667.1113
$ws.Cells.Item(73, 8).Value = 4755.5117
$ws.Cells.Item(73, 9).Value = 3767.96
$ws.Cells.Item(73, 10).Value = 6127.1113
$ws.Cells.Item(73, 11).Value = 3767.96
$ws.Cells.Item(73, 12).Value = 6127.1113
$ws.Cells.Item(73, 13).Value = -2831.96
$ws.Cells.Item(73, 14).Value = -7999.1113
$ws.Cells.Item(102, 8).Value = 3138.7942
$ws.Cells.Item(102, 9).Value = 1490.9131
$ws.Cells.Item(102, 10).Value = 6584.364
$ws.Cells.Item(102, 11).Value = 1490.9131
$ws.Cells.Item(102, 12).Value = 6584.364
$ws.Cells.Item(102, 13).Value = 131.0869
$ws.Cells.Item(102, 14).Value = -9828.364
$ws.Cells.Item(107, 8).Value = 5486.316
$ws.Cells.Item(107, 9).Value = 9237.727999999999
$ws.Cells.Item(107, 10).Value = 328.125
$ws.Cells.Item(107, 11).Value = 9237.727999999999
$ws.Cells.Item(107, 12).Value = 328.125
$ws.Cells.Item(107, 13).Value = -7317.727999999999
$ws.Cells.Item(107, 14).Value = -4168.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(43, 8).Value = 4760.5713
$ws.Cells.Item(43, 9).Value = 3862
$ws.Cells.Item(43, 10).Value = 5120
$ws.Cells.Item(43, 11).Value = 3862
$ws.Cells.Item(43, 12).Value = 5120
$ws.Cells.Item(43, 13).Value = -3669
$ws.Cells.Item(43, 14).Value = -5506
$ws.Cells.Item(61, 8).Value = 1080.5
$ws.Cells.Item(61, 9).Value = 696.6
$ws.Cells.Item(61, 11).Value = 696.6
$ws.Cells.Item(61, 13).Value = -494.6
$ws.Cells.Item(113, 8).Value = 1080.5
$ws.Cells.Item(113, 9).Value = 696.6
$ws.Cells.Item(113, 11).Value = 696.6
$ws.Cells.Item(113, 13).Value = 1473.4
$ws.Cells.Item(132, 8).Value = 10826.962
$ws.Cells.Item(132, 9).Value = 3392.5386
$ws.Cells.Item(132, 10).Value = 18261.385
$ws.Cells.Item(132, 11).Value = 10177.6158
$ws.Cells.Item(132, 12).Value = 54784.155
$ws.Cells.Item(132, 13).Value = -7647.6158
$ws.Cells.Item(132, 14).Value = -59844.155
$ws.Cells.Item(136, 8).Value = 3898
$ws.Cells.Item(136, 9).Value = 2342.4443
$ws.Cells.Item(136, 10).Value = 8205.691999999999
$ws.Cells.Item(136, 11).Value = 7027.3329
$ws.Cells.Item(136, 12).Value = 24617.076
$ws.Cells.Item(136, 13).Value = -4477.3329
$ws.Cells.Item(136, 14).Value = -29717.076
$ws.Cells.Item(140, 8).Value = 58809.668
$ws.Cells.Item(140, 10).Value = 58809.668
$ws.Cells.Item(140, 12).Value = 58809.668
$ws.Cells.Item(140, 14).Value = -69169.66800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 45501300
$ws.Cells.Item(122, 9).Value = 47667890
$ws.Cells.Item(122, 11).Value = 143003670
$ws.Cells.Item(122, 13).Value = -143001220
$ws.Cells.Item(126, 8).Value = 1432.1578
$ws.Cells.Item(126, 9).Value = 887.2727
$ws.Cells.Item(126, 10).Value = 2181.375
$ws.Cells.Item(126, 11).Value = 2661.8181
$ws.Cells.Item(126, 12).Value = 6544.125
$ws.Cells.Item(126, 13).Value = -191.8181
$ws.Cells.Item(126, 14).Value = -11484.125
$ws.Cells.Item(45, 8).Value = 3880
$ws.Cells.Item(45, 9).Value = 0
$ws.Cells.Item(45, 10).Value = 3880
$ws.Cells.Item(45, 11).Value = 0
$ws.Cells.Item(45, 12).Value = 3880
$ws.Cells.Item(45, 13).ClearContents()
$ws.Cells.Item(45, 14).Value = -4862
